$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date value shared between C2 and C3 (edit in place so both references track)
$ws.Range("C2").Value = "июл 11 2017"
$ws.Range("C3").Value = "июл 11 2017"

# Add destination column (index / sequence numbers) in column A
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
